# Update countries & provincias Spain
# Applies the data refresh described by the commit:
#  - bumps the "Datos actualizados" timestamp
#  - updates/re-ranks a handful of country rows whose case counts changed
#    (Armenia moves above Argelia/Chequia; Oman, Hungria, Letonia and
#    Estado de Palestina get refreshed counts)
#  - fixes a couple of mis-ordered rows (Belice/Santa Lucia and
#    Islas Virgenes Britanicas/Papua Nueva Guinea) so the country name
#    matches its data again

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 09:05"

# 2) Oman (row 55) - refreshed counts
$ws.Range("E55").Value = 7984
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 43

# 3) Armenia now overtakes Argelia and Chequia in total cases, so it is
#    re-inserted right after Bolivia (row 58) with new data, and Argelia /
#    Chequia shift down one row each, keeping their own (unchanged) data.
$ws.Range("A58").Value = "Armenia"
$ws.Range("B58").Value = 9282
$ws.Range("C58").Value = 355
$ws.Range("D58").Value = 3386
$ws.Range("E58").Value = 5765
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 4
$ws.Range("H58").Value = 131

$ws.Range("A59").Value = "Argelia"
$ws.Range("B59").Value = 9267
$ws.Range("C59").Value = 0
$ws.Range("D59").Value = 5549
$ws.Range("E59").Value = 3072
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 646

$ws.Range("A60").Value = "Chequia"
$ws.Range("B60").Value = 9230
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 6546
$ws.Range("E60").Value = 2365
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 319

# 4) Hungria (row 75) - refreshed counts
$ws.Range("B75").Value = 3876
$ws.Range("C75").Value = 9
$ws.Range("D75").Value = 2147
$ws.Range("E75").Value = 1203
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = 526

# 5) Letonia (row 114) - refreshed counts
$ws.Range("B114").Value = 1066
$ws.Range("C114").Value = 1
$ws.Range("E114").Value = 297

# 6) Estado de Palestina (row 142) - refreshed counts
$ws.Range("B142").Value = 448
$ws.Range("C142").Value = 1
$ws.Range("E142").Value = 77

# 7) Belice / Santa Lucia were swapped - fix names back to their data
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

# 8) Islas Virgenes Britanicas / Papua Nueva Guinea were swapped - fix names
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
